$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (column D) text value for rows whose price changed
$priceUpdates = @{
    2  = "244.39"
    3  = "24.01"
    4  = "5.255"
    5  = "0.05840"
    6  = "6.456"
    7  = "3.331"
    8  = "0.8076"
    9  = "0.8965"
    11 = "0.07096"
    12 = "0.03056"
    13 = "0.03028"
    14 = "0.09338"
    15 = "3.812"
    16 = "0.001542"
    17 = "0.04698"
    18 = "0.0006047"
    19 = "0.006195"
    20 = "0.001255"
    21 = "0.004074"
    22 = "0.00008695"
    24 = "2.169"
    25 = "0.3186"
    28 = "0.0002327"
    40 = "0.03843"
    41 = "0.006296"
    42 = "0.1054"
    43 = "0.002599"
    44 = "0.006973"
    45 = "0.00005303"
    47 = "0.4850"
    48 = "0.006157"
}

# Update the price cells, forcing text storage so the exact numeric-looking
# string (including trailing zeros / significant digits) is preserved.
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Column G (Hora) for rows 2..51 changes from "17" to "18" (also textual).
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = "@"
    $cell.Value = "18"
}
